$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected (legacy password "D382"); unprotect before
# writing, then re-apply the same protection so the sheet ends up in the
# same protected state it started in.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclosure banner (A16).
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-05 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) / Percent Change (E) columns for rows 2-12.
$ws.Range("D2").Value = 0.03102358496500318
$ws.Range("E2").Value = -0.01702417432754511

$ws.Range("D3").Value = 0.02396404680865104
$ws.Range("E3").Value = -0.01582867783985109

$ws.Range("D4").Value = 0.05148652526780813
$ws.Range("E4").Value = -0.001175640724194738

$ws.Range("D5").Value = 0.1378869120879147
$ws.Range("E5").Value = 0.0004880429477793768

$ws.Range("D6").Value = 0.0304060637897119
$ws.Range("E6").Value = 0.03430656934306575

$ws.Range("D7").Value = 0.1207202269458563
$ws.Range("E7").Value = -0.003995505056811233

$ws.Range("D8").Value = 0.1026314733796508
$ws.Range("E8").Value = 0

$ws.Range("D9").Value = 0.02896344032523652
$ws.Range("E9").Value = 0.0121212121212122

$ws.Range("D10").Value = 0.1249284622547748
$ws.Range("E10").Value = 0.008170245107353313

$ws.Range("D11").Value = 0.2449635172335052
$ws.Range("E11").Value = -0.002117278836417213

$ws.Range("D12").Value = 0.1030257469418875
$ws.Range("E12").Value = -0.004838397522740512

# Row 13 "Total" only has its Percent Change refreshed (Weight stays 1).
$ws.Range("E13").Value = 0.00001471701436384087

# Restore sheet protection with the original password.
$ws.Protect("D382")
